# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Offense (OFF) sheet — row 2 totals after the Wild Card round
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 415
$wsOff.Range("C2").Value = 286
$wsOff.Range("D2").Value = 60

# Defense (DEF) sheet — row 2 totals after the Wild Card round
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 525
$wsDef.Range("C2").Value = 350
$wsDef.Range("D2").Value = 144
$wsDef.Range("E2").Value = 65
$wsDef.Range("F2").Value = 12
